$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("test_storageOperations.py")

# --- test_storageOperations.py: make room for two new section header rows ---
# Shift the existing 22 data rows (old rows 1-11 and 13-23) down by 11 so the
# "Puts" table lands at rows 12-22 and the "Take" table at rows 24-34, leaving
# rows 1 and 11 free for new section headers.
$ws3.Rows("1:11").Insert()

# Row 11 header ("Storage Take Tests") is written first so it claims the
# earlier shared-string slot, matching the source workbook's string order.
$ws3.Range("A11").Value = "Storage Take Tests"
$ws3.Range("A11").Font.Bold = $true

# Row 1 header ("Storage Puts Tests") written second.
$ws3.Range("A1").Value = "Storage Puts Tests"
$ws3.Range("A1").Font.Bold = $true

# Fix the typo: the "puts into storage" test columns D and E (rows 33 and 34
# in the shifted layout) were wrongly using formulas copied from the "take"
# table; replace them with plain 0 literals.
$ws3.Range("E32").Value = 0
$ws3.Range("D33").Value = 0
$ws3.Range("E33").Value = 0

# New column E needs a width now that it holds real data.
$ws3.Columns.Item(5).ColumnWidth = 9.2

# Update the view: select A6, and make this the active/selected sheet
# (test_modelLogic.py was previously the selected tab, so activating this
# sheet also clears that sheet's tabSelected flag).
$ws3.Range("A6").Select() | Out-Null
$ws3.Activate()
